# Convert the M2Doc field codes (w:fldChar/w:instrText) used in this
# document into plain literal text runs wrapped in curly braces, e.g.
#   { begin }<instrText>m:link 'bookmark1'</instrText>{ end }
# becomes
#   {m:link 'bookmark1'}
#
# This mirrors the "TokenIteratorFieldRewriterSplit" rewrite: every Word
# field built from M2Doc instruction text is unlinked into ordinary text
# so the template text shows the raw "{m:...}" token instead of a Word
# field.

$d = $word.ActiveDocument

# Walk the Fields collection back-to-front so that deleting/inserting
# around one field never invalidates the character offsets of the
# fields that come before it in the document.
for ($idx = $d.Fields.Count; $idx -ge 1; $idx--) {

    $f = $d.Fields.Item($idx)
    $code = $f.Code.Text
    $trimmed = $code.Trim()

    # Position of the field's opening fldChar (one character before the
    # field instruction text begins).
    $fieldStart = $f.Code.Start - 1

    # The "m:bookmark self." field is preceded by a standalone run of
    # four spaces that the rewrite also removes (see the expected diff:
    # "Test bookmark :     {m:bookmark self.}" becomes
    # "Test bookmark : {m:bookmark self.}").
    if ($trimmed -eq "m:bookmark self.") {
        $preceding = $d.Range($fieldStart - 4, $fieldStart)
        if ($preceding.Text -eq "    ") {
            $preceding.Delete()
            $fieldStart = $fieldStart - 4
        }
    }

    # Remove the field (fldChar begin/end + instrText runs) entirely.
    $f.Delete()

    # Re-insert the field instruction as literal text, wrapped in the
    # curly braces M2Doc uses for its textual token syntax.
    $newText = "{" + $trimmed + "}"
    $insertionPoint = $d.Range($fieldStart, $fieldStart)
    $insertionPoint.InsertAfter($newText)
}
